$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Remove the Taxon_Local_ID / ${iNaturalistTaxonId} column and the
# suborder / infraorder / superfamily columns from the Materials sheet.
# Delete right-to-left so earlier column letters stay valid.
$ws.Columns("AT").Delete()
$ws.Columns("AS").Delete()
$ws.Columns("AR").Delete()
$ws.Columns("A").Delete()

# Fix the lingering "${summary.Author}" template placeholder -> "${summary.authority}"
$ws.Range("AX2").Value = '${summary.authority}'
